$d = $word.ActiveDocument

$replacements = @(
    @("340÷9=", "908÷7="),
    @("446÷8=", "886÷4="),
    @("107÷6=", "127÷7="),
    @("782÷5=", "616÷8="),
    @("777÷8=", "558÷4="),
    @("566÷6=", "149÷2="),
    @("558÷3=", "583÷7="),
    @("906÷7=", "994÷3="),
    @("762÷2=", "645÷8="),
    @("867÷7=", "489÷3="),
    @("657÷7=", "671÷5="),
    @("505÷6=", "281÷6="),
    @("168÷5=", "687÷7="),
    @("876÷7=", "913÷9="),
    @("331÷3=", "899÷6="),
    @("704÷9=", "614÷6="),
    @("415÷9=", "129÷2="),
    @("561÷8=", "334÷3="),
    @("332÷8=", "384÷7="),
    @("423÷3=", "164÷9="),
    @("468÷4=", "494÷5="),
    @("740÷4=", "598÷7="),
    @("409÷6=", "416÷4="),
    @("225÷7=", "137÷5="),
    @("745÷2=", "628÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
